$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 02:16"

# Row 4 - China
$ws.Range("B4").Value = 81054
$ws.Range("C4").Value = 46
$ws.Range("D4").Value = 72440
$ws.Range("E4").Value = 5353
$ws.Range("F4").Value = 1845
$ws.Range("G4").Value = 6
$ws.Range("H4").Value = 3261

# Row 6 - Estados Unidos
$ws.Range("B6").Value = 26112
$ws.Range("C6").Value = 6729
$ws.Range("G6").Value = 69
$ws.Range("H6").Value = 325

# Row 11 - Corea del Sur
$ws.Range("B11").Value = 8897
$ws.Range("C11").Value = 245
$ws.Range("D11").Value = 2909
$ws.Range("E11").Value = 5884
$ws.Range("G11").Value = 10
$ws.Range("H11").Value = 104
